$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "BACKBONE v15`ndefault value`nNULL`narg_check`n`"`"`nnot required sections`nwarning before output"
$ws.Range("C1").Value = "colons_check"
$ws.Range("D1").Value = "all_args_here"
$ws.Range("E1").Value = "safer_check`nlib_path`nerror_text`nmanual in good order"
$ws.Range("F1").Value = "81`ninternal errors number"
$ws.Range("G1").Value = "ERROR`ntempo.cat <-"
$ws.Range("H1").Value = "59`nerror message end by ."
$ws.Range("I1").Value = "71`nerror EMPTY"
$ws.Range("J1").Value = "76`nerror_text = embed"
$ws.Range("K1").Value = "77`nsub(pattern = `"^ERROR IN `""
$ws.Range("L1").Value = "78`ninternal fun no arg_check -> deal with `"`""
$ws.Range("M1").Value = "80`n. preceeded by :::`n\b(?!base\b)\w+(?=::)`ndo it after colons_check"
$ws.Range("N1").Value = "79`nrequired function writted in internal"
$ws.Range("O1").Value = "check lib_path = lib_path or not"
$ws.Range("P1").Value = "in main fun, check that intern fun have the good args"
$ws.Range("Q1").Value = "Unit tests`nbackbone v15.2"
$ws.Range("R1").Value = "    #### warning output`nuse the backbone one`nckeck no ini.warning.length`nwarn.count"
$ws.Range("S1").Value = "backbone v15.2: warning before output"
$ws.Range("T1").Value = "backbone v15.2: `nin internal functions`nerror_text argument with no value but comment that no tracable error message returned "
$ws.Range("U1").Value = "data1 argument -> data`nbecause protected by no arg value section"
$ws.Range("M2").Value = $null
$ws.Range("P2").Value = "x"
$ws.Range("R2").Value = "x"
$ws.Range("S2").Value = "x"
$ws.Range("T2").Value = "x"
$ws.Range("U2").Value = "x"
$ws.Range("M3").Value = $null
$ws.Range("Q3").Value = "x"
$ws.Range("R3").Value = "x"
$ws.Range("S3").Value = "x"
$ws.Range("T3").Value = "x"
$ws.Range("U3").Value = "x"
$ws.Range("D4").Value = "`"+check all arg_check`""
$ws.Range("O4").Value = "x"
$ws.Range("P4").Value = "x"
$ws.Range("R4").Value = "x"
$ws.Range("S4").Value = "x"
$ws.Range("T4").Value = "x"
$ws.Range("U4").Value = "x"
$ws.Range("D5").Value = $null
$ws.Range("M5").Value = $null
$ws.Range("P5").Value = "x"
$ws.Range("R5").Value = "x"
$ws.Range("S5").Value = "x"
$ws.Range("T5").Value = "x"
$ws.Range("U5").Value = "x"
$ws.Range("D6").Value = $null
$ws.Range("O6").Value = "x"
$ws.Range("P6").Value = "x"
$ws.Range("R6").Value = "x"
$ws.Range("S6").Value = "x"
$ws.Range("T6").Value = "x"
$ws.Range("U6").Value = "x"
$ws.Range("D7").Value = $null
$ws.Range("M7").Value = $null
$ws.Range("P7").Value = "x"
$ws.Range("R7").Value = "x"
$ws.Range("S7").Value = "x"
$ws.Range("T7").Value = "x"
$ws.Range("U7").Value = "x"
$ws.Range("D8").Value = $null
$ws.Range("M8").Value = $null
$ws.Range("P8").Value = "x"
$ws.Range("R8").Value = "x"
$ws.Range("S8").Value = "x"
$ws.Range("T8").Value = "x"
$ws.Range("U8").Value = "x"
$ws.Range("M9").Value = $null
$ws.Range("Q9").Value = "x"
$ws.Range("R9").Value = "x"
$ws.Range("S9").Value = "x"
$ws.Range("T9").Value = "x"
$ws.Range("U9").Value = "x"
$ws.Range("D10").Value = "issue 82"
$ws.Range("M10").Value = $null
$ws.Range("P10").Value = "x"
$ws.Range("R10").Value = "x"
$ws.Range("S10").Value = "x"
$ws.Range("T10").Value = "x"
$ws.Range("U10").Value = "x"
$ws.Range("M11").Value = $null
$ws.Range("P11").Value = "x"
$ws.Range("R11").Value = "x"
$ws.Range("S11").Value = "x"
$ws.Range("T11").Value = "x"
$ws.Range("U11").Value = "x"
$ws.Range("D12").Value = $null
$ws.Range("M12").Value = $null
$ws.Range("P12").Value = "x"
$ws.Range("R12").Value = "x"
$ws.Range("S12").Value = "x"
$ws.Range("T12").Value = "x"
$ws.Range("U12").Value = "x"
$ws.Range("D13").Value = $null
$ws.Range("M13").Value = $null
$ws.Range("P13").Value = "x"
$ws.Range("R13").Value = "x"
$ws.Range("S13").Value = "x"
$ws.Range("T13").Value = "x"
$ws.Range("U13").Value = "x"
$ws.Range("D14").Value = $null
$ws.Range("M14").Value = $null
$ws.Range("P14").Value = "x"
$ws.Range("R14").Value = "x"
$ws.Range("S14").Value = "x"
$ws.Range("T14").Value = "x"
$ws.Range("U14").Value = "x"
$ws.Range("D15").Value = $null
$ws.Range("M15").Value = $null
$ws.Range("P15").Value = "x"
$ws.Range("R15").Value = "x"
$ws.Range("S15").Value = "x"
$ws.Range("T15").Value = "x"
$ws.Range("U15").Value = "x"
$ws.Range("D16").Value = $null
$ws.Range("M16").Value = $null
$ws.Range("P16").Value = "x"
$ws.Range("R16").Value = "x"
$ws.Range("S16").Value = "x"
$ws.Range("T16").Value = "x"
$ws.Range("U16").Value = "x"
$ws.Range("D17").Value = $null
$ws.Range("M17").Value = $null
$ws.Range("P17").Value = "x"
$ws.Range("R17").Value = "x"
$ws.Range("S17").Value = "x"
$ws.Range("T17").Value = "x"
$ws.Range("U17").Value = "x"
$ws.Range("M18").Value = $null
$ws.Range("Q18").Value = "x"
$ws.Range("R18").Value = "x"
$ws.Range("S18").Value = "x"
$ws.Range("T18").Value = "x"
$ws.Range("U18").Value = "x"
$ws.Range("D19").Value = $null
$ws.Range("M19").Value = $null
$ws.Range("P19").Value = "x"
$ws.Range("R19").Value = "x"
$ws.Range("S19").Value = "x"
$ws.Range("T19").Value = "x"
$ws.Range("U19").Value = "x"
$ws.Range("D20").Value = $null
$ws.Range("O20").Value = "x"
$ws.Range("P20").Value = "x"
$ws.Range("R20").Value = "x"
$ws.Range("S20").Value = "x"
$ws.Range("T20").Value = "x"
$ws.Range("U20").Value = "x"
$ws.Range("D21").Value = $null
$ws.Range("O21").Value = "x"
$ws.Range("P21").Value = "x"
$ws.Range("R21").Value = "x"
$ws.Range("S21").Value = "x"
$ws.Range("T21").Value = "x"
$ws.Range("U21").Value = "x"
$ws.Range("D22").Value = $null
$ws.Range("O22").Value = "x"
$ws.Range("P22").Value = "x"
$ws.Range("R22").Value = "x"
$ws.Range("S22").Value = "x"
$ws.Range("T22").Value = "x"
$ws.Range("U22").Value = "x"
